$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table

$tbl.Cell(5, 1).Shape.TextFrame.TextRange.Text = "2018.07.27"
$tbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "게임판 크기 추가"
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = "햄과함께"
